$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (65) describing the "missing_weight_values" qc_flag
$ws.Range("A65").Value = "missing_weight_values"
$ws.Range("B65").Value = "Subjects"
$ws.Range("C65").Value = "Subject record missing weight value"
$ws.Range("D65").Value = "Soft Stop (Missing Preferred Column)"

# Reflect the cursor position left in the saved file (scrolled view + active cell)
$ws.Range("D65").Select()
